# Weekly data refresh: insert a new price record as row 20 (pushing the
# existing rows 20-117 down to 21-118), matching the "Fruta / hortaliza,
# semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a fresh row at position 20.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A20").Value = 11
$ws.Range("B20").Value = 'Vega Monumental Concepción'
$ws.Range("C20").Value = 'Bíobío'
$ws.Range("D20").Value = 44749
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = 'Fruta'
$ws.Range("G20").Value = 100108
$ws.Range("H20").Value = 'Tropicales y subtropicales'
$ws.Range("I20").Value = 100108002
$ws.Range("J20").Value = 'Mango'
$ws.Range("K20").Value = 'Sin especificar'
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 7500
$ws.Range("O20").Value = 8000
$ws.Range("P20").Value = 7750
$ws.Range("Q20").Value = '$/bandeja 4 kilos'
$ws.Range("R20").Value = 'Brasil'
$ws.Range("S20").Value = 1938
$ws.Range("T20").Value = 4
